# Update Data Model and XML file
#  - Rename the "Role" sheet to "Role List"
#  - Move the saved selection on the "Role List" sheet from E11 to C40
#  - Move the saved selection on the "License List" sheet from D10 to B22
#    (selected last so "License List" remains the active/visible tab)

$wb = $excel.ActiveWorkbook

$roleSheet = $wb.Worksheets.Item("Role")
$roleSheet.Name = "Role List"
$roleSheet.Range("C40").Select()

$licenseSheet = $wb.Worksheets.Item("License List")
$licenseSheet.Range("B22").Select()
